# Auto-generated Excel COM-interop edit script.
# Applies the 'cryptos' price/volume refresh described by the commit diff:
# "Updated cryptos list on Wed Dec  6 05:44:23 UTC 2023 with GitHub Actions"
#
# The sheet stores Price (D) / Volume(1h) (E) figures as literal text (they
# include things like thousands-dot-separated prices, '  +x.xx%  ' padding,
# and trailing zeros that must be preserved exactly). A handful of the new
# Price strings parse as plain numbers AND would lose a trailing zero if
# Excel auto-converted them (e.g. '1.20' -> 1.2, '0.000220' -> 0.00022), so
# those specific cells are pinned to text format before the value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.742.59'
$ws.Range("E2").Value = '  +4.28%  '
$ws.Range("D3").Value = '2.267.18'
$ws.Range("E3").Value = '  +1.48%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '229.56'
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '63.26'
$ws.Range("E7").Value = '  +4.11%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = '0.421'
$ws.Range("E9").Value = '  +3.47%  '
$ws.Range("D10").Value = '0.0979'
$ws.Range("E10").Value = '  +8.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.30'
$ws.Range("E11").Value = '  -1.60%  '
$ws.Range("D12").Value = '25.77'
$ws.Range("E12").Value = '  +13.16%  '
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").Value = '2.607.44'
$ws.Range("E14").Value = '  +1.62%  '
$ws.Range("D15").Value = '15.55'
$ws.Range("E15").Value = '  -0.01%  '
$ws.Range("D16").Value = '5.84'
$ws.Range("E16").Value = '  +3.05%  '
$ws.Range("D17").Value = '0.809'
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").Value = '2.275.61'
$ws.Range("E18").Value = '  +1.64%  '
$ws.Range("D19").Value = '43.629.40'
$ws.Range("E19").Value = '  +4.19%  '
$ws.Range("D20").Value = '0.0₃0942'
$ws.Range("E20").Value = '  +3.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.80'
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").Value = '6.06'
$ws.Range("E22").Value = '  -2.12%  '
$ws.Range("D23").Value = '247.95'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").Value = '2.48'
$ws.Range("E25").Value = '  +3.64%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").Value = '9.92'
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("D28").Value = '170.87'
$ws.Range("E28").Value = '  +0.69%  '
$ws.Range("D29").Value = '0.137'
$ws.Range("E29").Value = '  -3.52%  '
$ws.Range("D30").Value = '20.44'
$ws.Range("E30").Value = '  +2.42%  '
$ws.Range("E31").Value = '  +2.52%  '
$ws.Range("E32").Value = '  +10.05%  '
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("D34").Value = '0.0687'
$ws.Range("E34").Value = '  +4.05%  '
$ws.Range("D35").Value = '5.07'
$ws.Range("E35").Value = '  +0.62%  '
$ws.Range("D36").Value = '4.65'
$ws.Range("E36").Value = '  -1.38%  '
$ws.Range("D37").Value = '6.72'
$ws.Range("E37").Value = '  +1.72%  '
$ws.Range("D38").Value = '3.75'
$ws.Range("E38").Value = '  +3.25%  '
$ws.Range("D39").Value = '2.31'
$ws.Range("E39").Value = '  -3.94%  '
$ws.Range("D40").Value = '0.0245'
$ws.Range("E40").Value = '  +1.57%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.000220'
$ws.Range("E42").Value = '  -10.34%  '
$ws.Range("D43").Value = '8.32'
$ws.Range("E43").Value = '  -4.16%  '
$ws.Range("D44").Value = '10.44'
$ws.Range("E44").Value = '  +18.98%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = '0.0961'
$ws.Range("E45").Value = '  -0.42%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").Value = '4.44'
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.20'
$ws.Range("E47").Value = '  -2.61%  '
$ws.Range("D48").Value = '96.59'
$ws.Range("E48").Value = '  -2.55%  '
$ws.Range("D49").Value = '1.471.39'
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("D50").Value = '16.69'
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").Value = '2.28'
$ws.Range("E51").Value = '  -0.09%  '
